# Auto-generated edit script: updates currentAveragePrice-derived columns
# (H/I/J/K/L/M/N) across multiple job sheets, matching the scheduled-runner
# price refresh captured in the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 724.4286
$ws.Range("I6").Value = 724.4286
$ws.Range("K6").Value = 2173.2858
$ws.Range("M6").Value = -2061.2858
$ws.Range("H28").Value = 827.6923
$ws.Range("I28").Value = 861.6667
$ws.Range("K28").Value = 861.6667
$ws.Range("M28").Value = -376.6667
$ws.Range("H98").Value = 2286.5789
$ws.Range("I98").Value = 2279.353
$ws.Range("K98").Value = 2279.353
$ws.Range("M98").Value = -781.3530000000001
$ws.Range("H106").Value = 250003200
$ws.Range("I106").Value = 250003200
$ws.Range("K106").Value = 250003200
$ws.Range("M106").Value = -250002569
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H122").Value = 2286.5789
$ws.Range("I122").Value = 2279.353
$ws.Range("K122").Value = 6838.059
$ws.Range("M122").Value = -4388.059
$ws.Range("H132").Value = 27720.4
$ws.Range("I132").Value = 41507.46
$ws.Range("K132").Value = 124522.38
$ws.Range("M132").Value = -121992.38
$ws.Range("H138").Value = 9323.885
$ws.Range("I138").Value = 13951.286
$ws.Range("J138").Value = 3925.25
$ws.Range("K138").Value = 41853.858
$ws.Range("L138").Value = 11775.75
$ws.Range("M138").Value = -36713.858
$ws.Range("N138").Value = -22055.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 146768.58
$ws.Range("I45").Value = 146768.58
$ws.Range("K45").Value = 146768.58
$ws.Range("M45").Value = -146391.58
$ws.Range("H61").Value = 7579.6787
$ws.Range("I61").Value = 7265
$ws.Range("K61").Value = 7265
$ws.Range("M61").Value = -7053
$ws.Range("H102").Value = 2845.3333
$ws.Range("I102").Value = 2836.0667
$ws.Range("J102").Value = 2891.6667
$ws.Range("K102").Value = 2836.0667
$ws.Range("L102").Value = 2891.6667
$ws.Range("M102").Value = -1214.0667
$ws.Range("N102").Value = -6135.6667
$ws.Range("H122").Value = 3121
$ws.Range("I122").Value = 2950.5454
$ws.Range("K122").Value = 8851.636200000001
$ws.Range("M122").Value = -6401.636200000001
$ws.Range("H125").Value = 74995
$ws.Range("J125").Value = 74995
$ws.Range("L125").Value = 74995
$ws.Range("N125").Value = -84835
$ws.Range("H132").Value = 2285.3333
$ws.Range("I132").Value = 1467.1333
$ws.Range("J132").Value = 4330.8335
$ws.Range("K132").Value = 4401.3999
$ws.Range("L132").Value = 12992.5005
$ws.Range("M132").Value = -1871.3999
$ws.Range("N132").Value = -18052.5005
$ws.Range("H136").Value = 7579.6787
$ws.Range("I136").Value = 7265
$ws.Range("K136").Value = 21795
$ws.Range("M136").Value = -19245

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 8610.65
$ws.Range("I105").Value = 10033.692
$ws.Range("K105").Value = 10033.692
$ws.Range("M105").Value = -8286.691999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2157.5
$ws.Range("I58").Value = 1627
$ws.Range("J58").Value = 2591.5454
$ws.Range("K58").Value = 1627
$ws.Range("L58").Value = 2591.5454
$ws.Range("M58").Value = -1424
$ws.Range("N58").Value = -2997.5454
$ws.Range("H86").Value = 15649.762
$ws.Range("J86").Value = 23057.3
$ws.Range("L86").Value = 23057.3
$ws.Range("N86").Value = -25303.3
$ws.Range("H89").Value = 15649.762
$ws.Range("J89").Value = 23057.3
$ws.Range("L89").Value = 115286.5
$ws.Range("N89").Value = -126518.5
$ws.Range("H105").Value = 1748.4286
$ws.Range("I105").Value = 1252.5
$ws.Range("K105").Value = 1252.5
$ws.Range("M105").Value = 494.5
$ws.Range("H107").Value = 1486.3667
$ws.Range("I107").Value = 1474.2941
$ws.Range("K107").Value = 1474.2941
$ws.Range("M107").Value = 445.7058999999999
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 13336536
$ws.Range("I132").Value = 3106.7334
$ws.Range("J132").Value = 33336680
$ws.Range("K132").Value = 9320.200199999999
$ws.Range("L132").Value = 100010040
$ws.Range("M132").Value = -6790.200199999999
$ws.Range("N132").Value = -100015100
$ws.Range("H134").Value = 2025.3793
$ws.Range("I134").Value = 1481.9412
$ws.Range("J134").Value = 2795.25
$ws.Range("K134").Value = 4445.8236
$ws.Range("L134").Value = 8385.75
$ws.Range("M134").Value = -1910.8236
$ws.Range("N134").Value = -13455.75
$ws.Range("H136").Value = 2157.5
$ws.Range("I136").Value = 1627
$ws.Range("J136").Value = 2591.5454
$ws.Range("K136").Value = 4881
$ws.Range("L136").Value = 7774.6362
$ws.Range("M136").Value = -2331
$ws.Range("N136").Value = -12874.6362
$ws.Range("H141").Value = 460988.1
$ws.Range("J141").Value = 406653.56
$ws.Range("L141").Value = 406653.56
$ws.Range("N141").Value = -417013.56

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 534.7143
$ws.Range("I5").Value = 524
$ws.Range("K5").Value = 1572
$ws.Range("M5").Value = -1460
$ws.Range("H7").Value = 297.25
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H92").Value = 439.06668
$ws.Range("I92").Value = 385.57144
$ws.Range("J92").Value = 485.875
$ws.Range("K92").Value = 1156.71432
$ws.Range("L92").Value = 1457.625
$ws.Range("M92").Value = 91.28567999999996
$ws.Range("N92").Value = -3953.625
$ws.Range("I107").Value = 754.4
$ws.Range("J107").Value = 52631804
$ws.Range("K107").Value = 2263.2
$ws.Range("L107").Value = 157895412
$ws.Range("M107").Value = -343.1999999999998
$ws.Range("N107").Value = -157899252
$ws.Range("H133").Value = 3570.9333
$ws.Range("I133").Value = 3468.8572
$ws.Range("K133").Value = 10406.5716
$ws.Range("M133").Value = -5346.571599999999
$ws.Range("H135").Value = 534.7143
$ws.Range("I135").Value = 524
$ws.Range("K135").Value = 4716
$ws.Range("M135").Value = -2181
$ws.Range("H140").Value = 14707899
$ws.Range("I140").Value = 25002128
$ws.Range("K140").Value = 75006384
$ws.Range("M140").Value = -75001204

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13293.833
$ws.Range("I70").Value = 17186.25
$ws.Range("J70").Value = 10179.9
$ws.Range("K70").Value = 17186.25
$ws.Range("L70").Value = 10179.9
$ws.Range("M70").Value = -16916.25
$ws.Range("N70").Value = -10719.9
$ws.Range("H73").Value = 13293.833
$ws.Range("I73").Value = 17186.25
$ws.Range("J73").Value = 10179.9
$ws.Range("K73").Value = 17186.25
$ws.Range("L73").Value = 10179.9
$ws.Range("M73").Value = -16250.25
$ws.Range("N73").Value = -12051.9
$ws.Range("H132").Value = 10004183
$ws.Range("I132").Value = 3721.9
$ws.Range("K132").Value = 11165.7
$ws.Range("M132").Value = -8635.700000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 965.61536
$ws.Range("I16").Value = 553
$ws.Range("K16").Value = 553
$ws.Range("M16").Value = -383
$ws.Range("H61").Value = 5558023
$ws.Range("H100").Value = 3000
$ws.Range("I100").Value = 3000
$ws.Range("K100").Value = 3000
$ws.Range("M100").Value = -2459
$ws.Range("H113").Value = 5558023
$ws.Range("H136").Value = 4712.625
$ws.Range("I136").Value = 4624.375
$ws.Range("J136").Value = 4800.875
$ws.Range("K136").Value = 13873.125
$ws.Range("L136").Value = 14402.625
$ws.Range("M136").Value = -11323.125
$ws.Range("N136").Value = -19502.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 9999
$ws.Range("J39").Value = 9999
$ws.Range("L39").Value = 9999
$ws.Range("N39").Value = -10825
$ws.Range("H132").Value = 61134.06
$ws.Range("I132").Value = 145290.72
$ws.Range("K132").Value = 435872.16
$ws.Range("M132").Value = -433342.16
